# predictions updated using weather data
# Append 4 new weekly-prediction rows (56-59) for the prediction made on
# 2021-01-09, matching the existing table layout (KNN model row + 3
# bare-prediction rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the "day the prediction was made" and is always the same
# shared string already used by rows 52-55 ("2021-01-09"). Copying the
# existing cell instead of assigning a literal string keeps it as a plain
# shared-string cell (no style / no auto date-conversion), just like the
# rest of the column.
$ws.Range("A55").Copy($ws.Range("A56"))
$ws.Range("A55").Copy($ws.Range("A57"))
$ws.Range("A55").Copy($ws.Range("A58"))
$ws.Range("A55").Copy($ws.Range("A59"))

# Row 56 - week of 10 Jan -- 16 Jan 2021 (full KNN metrics row)
$ws.Range("B56").Value = "10 Jan -- 16 Jan 2021"
$ws.Range("C56").Value = 3333.57
$ws.Range("D56").Value = 1604.21
$ws.Range("E56").Value = 1729.36
$ws.Range("F56").Value = "KNN"
$ws.Range("J56").Value = 1495.63
$ws.Range("K56").Value = 48.09

# Row 57 - week of 17 Jan -- 23 Jan 2021
$ws.Range("B57").Value = "17 Jan -- 23 Jan 2021"
$ws.Range("D57").Value = 1915.89
$ws.Range("F57").Value = "KNN"

# Row 58 - week of 24 Jan -- 30 Jan 2021
$ws.Range("B58").Value = "24 Jan -- 30 Jan 2021"
$ws.Range("D58").Value = 2104.94
$ws.Range("F58").Value = "KNN"

# Row 59 - week of 31 Jan -- 06 Feb 2021
$ws.Range("B59").Value = "31 Jan -- 06 Feb 2021"
$ws.Range("D59").Value = 2208.2
$ws.Range("F59").Value = "KNN"
